# "excel utility implemented in login page"
# - rename Sheet1 -> PROJECTDATA
# - add a new sheet USERCREDENTIAL (right after PROJECTDATA) holding a
#   URL / username / password table, with the URL and password cells
#   turned into hyperlinks
# - leave the new sheet as the active/selected tab

$wb = $excel.ActiveWorkbook

# --- rename the existing sheet -------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "PROJECTDATA"
[void]$ws1.Range("E17").Select()

# --- add the credentials sheet right after PROJECTDATA --------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "USERCREDENTIAL"

$ws2.Range("A1").Value = "URL"
$ws2.Range("B1").Value = "username"
$ws2.Range("C1").Value = "password"

# data entry order matters for shared-string ordering: password, then
# username, then URL
$ws2.Range("C2").Value = "rmgy@9999"
$ws2.Range("B2").Value = "rmgyantra"
$ws2.Range("A2").Value = "http://localhost:8084/"

# turn the URL and password cells into hyperlinks
$ws2.Hyperlinks.Add($ws2.Range("C2"), "mailto:rmgy@9999")
$ws2.Hyperlinks.Add($ws2.Range("A2"), "http://localhost:8084/")

[void]$ws2.Range("H6").Select()
[void]$ws2.Activate()
